$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values (B1:E1)
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 data values (B2:E2)
$ws.Range("B2").Value = 10.233333326666665
$ws.Range("C2").Value = 4.8666666599999999
$ws.Range("D2").Value = 8.8333333266666667
$ws.Range("E2").Value = 10.899999993333333

# Row 3 data values (B3:E3)
$ws.Range("B3").Value = 4.5666666600000001
$ws.Range("C3").Value = 11.36666666
$ws.Range("D3").Value = 15.86666666
$ws.Range("E3").Value = 9.2333333266666653

# Update the sheet selection to match the new authored range
$ws.Range("B1:E3").Select()
